$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.773.05"
$ws.Range("E2").Value = "  +1.18%  "
$ws.Range("D3").Value = "3.157.56"
$ws.Range("E3").Value = "  +0.86%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "532.32"
$ws.Range("E5").Value = "  +0.32%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "140.11"
$ws.Range("E6").Value = "  +0.98%  "
$ws.Range("E7").Value = "  +0.02%  "
$ws.Range("E8").Value = "  +14.98%  "
$ws.Range("E9").Value = "  +0.45%  "
$ws.Range("E10").Value = "  +5.10%  "
$ws.Range("E11").Value = "  +2.64%  "
$ws.Range("E12").Value = "  +2.75%  "
$ws.Range("D13").Value = "3.701.94"
$ws.Range("E13").Value = "  +0.85%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "25.96"
$ws.Range("E14").Value = "  +1.59%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000171"
$ws.Range("E15").Value = "  +4.72%  "
$ws.Range("D16").Value = "58.813.74"
$ws.Range("E16").Value = "  +1.15%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.24"
$ws.Range("E17").Value = "  +4.07%  "
$ws.Range("D18").Value = "3.149.24"
$ws.Range("E18").Value = "  +0.81%  "
$ws.Range("E19").Value = "  +2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "8.20"
$ws.Range("E20").Value = "  +1.44%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "370.34"
$ws.Range("E21").Value = "  +4.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.81"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "70.01"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("E25").Value = "  +2.84%  "
$ws.Range("E26").Value = "  +0.13%  "
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.18"
$ws.Range("E28").Value = "  +12.26%  "
$ws.Range("E29").Value = "  -2.61%  "
$ws.Range("E30").Value = "  +1.34%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.10"
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "21.98"
$ws.Range("E32").Value = "  +3.11%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "5.20"
$ws.Range("E33").Value = "  +4.29%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.16"
$ws.Range("E34").Value = "  +0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "159.35"
$ws.Range("E35").Value = "  +0.47%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.28"
$ws.Range("E36").Value = "  +3.41%  "
$ws.Range("E37").Value = "  +6.98%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "25.30"
$ws.Range("E38").Value = "  -3.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.68"
$ws.Range("E39").Value = "  -0.07%  "
$ws.Range("D40").Value = "2.642.02"
$ws.Range("E40").Value = "  +10.36%  "
$ws.Range("E41").Value = "  +1.78%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "4.24"
$ws.Range("E42").Value = "  +5.92%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "38.77"
$ws.Range("E43").Value = "  +3.25%  "
$ws.Range("E44").Value = "  +1.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0287"
$ws.Range("E45").Value = "  +7.32%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.06%  "
$ws.Range("D47").Value = "3.198.98"
$ws.Range("E47").Value = "  +1.16%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.985"
$ws.Range("E49").Value = "  +0.72%  "
$ws.Range("E50").Value = "  +2.73%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "20.23"
$ws.Range("E51").Value = "  +2.03%  "
